# Generate Report for Archive
#
# 1) Status text changes from "Ready for handoff" to "In Translation"
#    (appears on the Overview sheet in columns E/F row 2, and on the
#    per-language sheets in column C row 2).
# 2) The "Status" column is narrowed (was 17.2159881591797 character-width
#    units, now 13.4101845877511) on the Overview sheet (columns E & F)
#    and on each per-language sheet (column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# --- 1) Update the status text wherever it appears -------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2) Narrow the Status column(s) -----------------------------------------
# ColumnWidth is expressed in "characters" and Excel snaps it to whole
# pixels, so we target the input value whose rounded result lands closest
# to the desired stored width (13.4101845877511).
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C
